$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expansion List")

# Copy the style/formatting of the last existing data row (row 13) down into
# the three new rows (14-16) so the new rows inherit the same cell style (s="3").
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 14: 10901-7 / Display for 2021-09 / SNOMEDCT / 2021-09 / 2.16.840.1.113883.6.96 / FN
$ws.Range("B14").Value = "Display for 2021-09"
$ws.Range("D14").Value = "2021-09"
$ws.Range("A14").Value = "10901-7"
$ws.Range("C14").Value = "SNOMEDCT"
$ws.Range("E14").Value = "2.16.840.1.113883.6.96"
$ws.Range("F14").Value = "FN"

# Row 15: 10901-8 / Display for 10901-8 / SNOMEDCT / 2021-09 / 2.16.840.1.113883.6.96 / FN
$ws.Range("B15").Value = "Display for 10901-8"
$ws.Range("A15").Value = "10901-8"
$ws.Range("C15").Value = "SNOMEDCT"
$ws.Range("D15").Value = "2021-09"
$ws.Range("E15").Value = "2.16.840.1.113883.6.96"
$ws.Range("F15").Value = "FN"

# Row 16: 10901-8 / Display for 10901-8 LOINC / LOINC / 2021-09 / 2.16.840.1.113883.6.1 / FN
$ws.Range("B16").Value = "Display for 10901-8 LOINC"
$ws.Range("C16").Value = "LOINC"
$ws.Range("E16").Value = "2.16.840.1.113883.6.1"
$ws.Range("A16").Value = "10901-8"
$ws.Range("D16").Value = "2021-09"
$ws.Range("F16").Value = "FN"

# Make "Expansion List" the active (tab-selected) sheet, with the newly added
# row selected, matching the saved selection/view state.
$ws.Activate()
$ws.Rows("16:16").Select()

$wb.Save()
